$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("first_eval")

$ws.Range("B2").Value = 0.4912602871908769
$ws.Range("C2").Value = 2.174891296644351
$ws.Range("D2").Value = 11.4276319550312
$ws.Range("E2").Value = 3.380478066047937
$ws.Range("F2").Value = 3.385133934666772
$ws.Range("G2").Value = 42

$ws.Range("B3").Value = 0.09227350360579856
$ws.Range("C3").Value = 2.692498561492474
$ws.Range("D3").Value = 10.01612202341713
$ws.Range("E3").Value = 3.164825749297602
$ws.Range("F3").Value = 3.202779701220204
$ws.Range("G3").Value = 41

$ws.Range("B4").Value = 0.5507705051130964
$ws.Range("C4").Value = 2.023147763338789
$ws.Range("D4").Value = 7.167790521971241
$ws.Range("E4").Value = 2.677272963664938
$ws.Range("F4").Value = 2.653385331149181
$ws.Range("G4").Value = 40

$ws.Range("B5").Value = 0.2467834876975515
$ws.Range("C5").Value = 2.086588795306851
$ws.Range("D5").Value = 9.964469520185652
$ws.Range("E5").Value = 3.156654799021529
$ws.Range("F5").Value = 3.188132309977417
$ws.Range("G5").Value = 39

$ws.Range("B6").Value = 0.5624371311069387
$ws.Range("C6").Value = 2.517633049209763
$ws.Range("D6").Value = 10.46235797449542
$ws.Range("E6").Value = 3.234556843602447
$ws.Range("F6").Value = 3.228039539870573
$ws.Range("G6").Value = 38

$ws.Range("B7").Value = 0.2885307811909213
$ws.Range("C7").Value = 2.261310828266383
$ws.Range("D7").Value = 7.061485019889288
$ws.Range("E7").Value = 2.657345483727942
$ws.Range("F7").Value = 2.678072976887184
$ws.Range("G7").Value = 37

$ws.Range("B8").Value = 0.543176441055371
$ws.Range("C8").Value = 2.381006566340435
$ws.Range("D8").Value = 9.179169249568099
$ws.Range("E8").Value = 3.029714384157045
$ws.Range("F8").Value = 3.022906027196908
$ws.Range("G8").Value = 36

$ws.Range("B9").Value = 0.1607468961940973
$ws.Range("C9").Value = 2.287366328590652
$ws.Range("D9").Value = 10.73307153727547
$ws.Range("E9").Value = 3.276136678662151
$ws.Range("F9").Value = 3.319962433532947
$ws.Range("G9").Value = 35

$ws.Range("B10").Value = 0.5939398229636811
$ws.Range("C10").Value = 2.070565034018146
$ws.Range("D10").Value = 7.383494754115199
$ws.Range("E10").Value = 2.717258683694874
$ws.Range("F10").Value = 2.691427627180946
$ws.Range("G10").Value = 34

$ws.Range("B11").Value = 0.3228496717043242
$ws.Range("C11").Value = 2.395276021476864
$ws.Range("D11").Value = 8.574011187212726
$ws.Range("E11").Value = 2.928141251239893
$ws.Range("F11").Value = 2.955411964361279
$ws.Range("G11").Value = 33

